# Regenerate s_val data to filter save games.
# Updates the numeric values in columns B:G for rows 2-5 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (2022-06-01)
$ws.Range("B2").Value = 0.04271373187048222
$ws.Range("C2").Value = 0.306821227259698
$ws.Range("D2").Value = 0.1494219747398047
$ws.Range("E2").Value = 0.4942365360607697
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.9931934699307545

# Row 3 (2022-05-27)
$ws.Range("B3").Value = 0.04271373187048222
$ws.Range("C3").Value = 1.655778082260271
$ws.Range("D3").Value = 3.537761648806719
$ws.Range("E3").Value = 10.19245300693656
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 15.42870646987403

# Row 4 (2022-05-05)
$ws.Range("B4").Value = 0.6606524410359556
$ws.Range("C4").Value = 1.655778082260271
$ws.Range("D4").Value = 3.537761648806719
$ws.Range("E4").Value = 0.4942365360607697
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 6.348428708163715

# Row 5 (2022-04-23)
$ws.Range("B5").Value = 3.286832544864788
$ws.Range("C5").Value = 1.655778082260271
$ws.Range("D5").Value = 0.7527432677738641
$ws.Range("E5").Value = 0.4942365360607697
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6.189590430959694

$wb.Save()
